# Increase the number of pre-allocated block rows in the results sheet
# (rows 2-5 already exist as placeholder rows of zeros; extend that same
# pattern down through row 17, i.e. add 12 more blocks).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 6; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
